$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> FAPs
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "App"
$ws.Cells.Item(2,3).Value = "Fpr2"
$ws.Cells.Item(2,4).Value = "FAPs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 103.4275383333333
$ws.Cells.Item(2,8).Value = 310.282615
$ws.Cells.Item(2,9).Value = 0.2485530285127421
$ws.Cells.Item(2,10).Value = 0.2485530285127421
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 1.257727
$ws.Cells.Item(2,14).Value = 3.773181
$ws.Cells.Item(2,15).Value = 0.3633293041311343
$ws.Cells.Item(2,16).Value = 0.3633293041311343
$ws.Cells.Item(2,17).Value = 130.0836075053683
$ws.Cells.Item(2,18).Value = 1170.752467548315
$ws.Cells.Item(2,19).Value = 0.09030659888922057
$ws.Cells.Item(2,20).Value = 0.09030659888922057

# Row 3: ECs -> M2
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "App"
$ws.Cells.Item(3,3).Value = "Fpr2"
$ws.Cells.Item(3,4).Value = "M2"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 103.4275383333333
$ws.Cells.Item(3,8).Value = 310.282615
$ws.Cells.Item(3,9).Value = 0.2485530285127421
$ws.Cells.Item(3,10).Value = 0.2485530285127421
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 2.112352666666667
$ws.Cells.Item(3,14).Value = 6.337058
$ws.Cells.Item(3,15).Value = 0.6102116154455982
$ws.Cells.Item(3,16).Value = 0.6102116154455982
$ws.Cells.Item(3,17).Value = 218.4754364051856
$ws.Cells.Item(3,18).Value = 1966.27892764667
$ws.Cells.Item(3,19).Value = 0.1516699450526562
$ws.Cells.Item(3,20).Value = 0.1516699450526562

# Row 4: ECs -> sCs
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "App"
$ws.Cells.Item(4,3).Value = "Fpr2"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 103.4275383333333
$ws.Cells.Item(4,8).Value = 310.282615
$ws.Cells.Item(4,9).Value = 0.2485530285127421
$ws.Cells.Item(4,10).Value = 0.2485530285127421
$ws.Cells.Item(4,11).Value = 1
$ws.Cells.Item(4,12).Value = 0.3333333333333333
$ws.Cells.Item(4,13).Value = 0.09159266666666667
$ws.Cells.Item(4,14).Value = 0.274778
$ws.Cells.Item(4,15).Value = 0.02645908042326748
$ws.Cells.Item(4,16).Value = 0.02645908042326749
$ws.Cells.Item(4,17).Value = 9.47320404271889
$ws.Cells.Item(4,18).Value = 85.25883638447002
$ws.Cells.Item(4,19).Value = 0.006576484570865339
$ws.Cells.Item(4,20).Value = 0.006576484570865339

# Row 5: FAPs -> FAPs
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "App"
$ws.Cells.Item(5,3).Value = "Fpr2"
$ws.Cells.Item(5,4).Value = "FAPs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 216.130539
$ws.Cells.Item(5,8).Value = 648.391617
$ws.Cells.Item(5,9).Value = 0.5193964865470273
$ws.Cells.Item(5,10).Value = 0.5193964865470272
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 1.257727
$ws.Cells.Item(5,14).Value = 3.773181
$ws.Cells.Item(5,15).Value = 0.3633293041311343
$ws.Cells.Item(5,16).Value = 0.3633293041311343
$ws.Cells.Item(5,17).Value = 271.833214424853
$ws.Cells.Item(5,18).Value = 2446.498929823677
$ws.Cells.Item(5,19).Value = 0.1887119640252875
$ws.Cells.Item(5,20).Value = 0.1887119640252875

# Row 6: FAPs -> M2
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "App"
$ws.Cells.Item(6,3).Value = "Fpr2"
$ws.Cells.Item(6,4).Value = "M2"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 216.130539
$ws.Cells.Item(6,8).Value = 648.391617
$ws.Cells.Item(6,9).Value = 0.5193964865470273
$ws.Cells.Item(6,10).Value = 0.5193964865470272
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 2.112352666666667
$ws.Cells.Item(6,14).Value = 6.337058
$ws.Cells.Item(6,15).Value = 0.6102116154455982
$ws.Cells.Item(6,16).Value = 0.6102116154455982
$ws.Cells.Item(6,17).Value = 456.543920404754
$ws.Cells.Item(6,18).Value = 4108.895283642786
$ws.Cells.Item(6,19).Value = 0.3169417691126294
$ws.Cells.Item(6,20).Value = 0.3169417691126294

# Row 7: FAPs -> sCs
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "App"
$ws.Cells.Item(7,3).Value = "Fpr2"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 216.130539
$ws.Cells.Item(7,8).Value = 648.391617
$ws.Cells.Item(7,9).Value = 0.5193964865470273
$ws.Cells.Item(7,10).Value = 0.5193964865470272
$ws.Cells.Item(7,11).Value = 1
$ws.Cells.Item(7,12).Value = 0.3333333333333333
$ws.Cells.Item(7,13).Value = 0.09159266666666667
$ws.Cells.Item(7,14).Value = 0.274778
$ws.Cells.Item(7,15).Value = 0.02645908042326748
$ws.Cells.Item(7,16).Value = 0.02645908042326749
$ws.Cells.Item(7,17).Value = 19.795972415114
$ws.Cells.Item(7,18).Value = 178.163751736026
$ws.Cells.Item(7,19).Value = 0.01374275340911036
$ws.Cells.Item(7,20).Value = 0.01374275340911036

# Row 8: M2 -> FAPs
$ws.Cells.Item(8,1).Value = "M2"
$ws.Cells.Item(8,2).Value = "App"
$ws.Cells.Item(8,3).Value = "Fpr2"
$ws.Cells.Item(8,4).Value = "FAPs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 71.607325
$ws.Cells.Item(8,8).Value = 214.821975
$ws.Cells.Item(8,9).Value = 0.1720839321833696
$ws.Cells.Item(8,10).Value = 0.1720839321833696
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 1.257727
$ws.Cells.Item(8,14).Value = 3.773181
$ws.Cells.Item(8,15).Value = 0.3633293041311343
$ws.Cells.Item(8,16).Value = 0.3633293041311343
$ws.Cells.Item(8,17).Value = 90.062466050275
$ws.Cells.Item(8,18).Value = 810.5621944524751
$ws.Cells.Item(8,19).Value = 0.06252313533233297
$ws.Cells.Item(8,20).Value = 0.06252313533233297

# Row 9: M2 -> M2
$ws.Cells.Item(9,1).Value = "M2"
$ws.Cells.Item(9,2).Value = "App"
$ws.Cells.Item(9,3).Value = "Fpr2"
$ws.Cells.Item(9,4).Value = "M2"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 71.607325
$ws.Cells.Item(9,8).Value = 214.821975
$ws.Cells.Item(9,9).Value = 0.1720839321833696
$ws.Cells.Item(9,10).Value = 0.1720839321833696
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 2.112352666666667
$ws.Cells.Item(9,14).Value = 6.337058
$ws.Cells.Item(9,15).Value = 0.6102116154455982
$ws.Cells.Item(9,16).Value = 0.6102116154455982
$ws.Cells.Item(9,17).Value = 151.2599239166167
$ws.Cells.Item(9,18).Value = 1361.33931524955
$ws.Cells.Item(9,19).Value = 0.1050076142498447
$ws.Cells.Item(9,20).Value = 0.1050076142498447

# Row 10: M2 -> sCs
$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "App"
$ws.Cells.Item(10,3).Value = "Fpr2"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 71.607325
$ws.Cells.Item(10,8).Value = 214.821975
$ws.Cells.Item(10,9).Value = 0.1720839321833696
$ws.Cells.Item(10,10).Value = 0.1720839321833696
$ws.Cells.Item(10,11).Value = 1
$ws.Cells.Item(10,12).Value = 0.3333333333333333
$ws.Cells.Item(10,13).Value = 0.09159266666666667
$ws.Cells.Item(10,14).Value = 0.274778
$ws.Cells.Item(10,15).Value = 0.02645908042326748
$ws.Cells.Item(10,16).Value = 0.02645908042326749
$ws.Cells.Item(10,17).Value = 6.558705849616667
$ws.Cells.Item(10,18).Value = 59.02835264655
$ws.Cells.Item(10,19).Value = 0.004553182601191883
$ws.Cells.Item(10,20).Value = 0.004553182601191884

# Row 11: sCs -> FAPs
$ws.Cells.Item(11,1).Value = "sCs"
$ws.Cells.Item(11,2).Value = "App"
$ws.Cells.Item(11,3).Value = "Fpr2"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 24.953198
$ws.Cells.Item(11,8).Value = 74.859594
$ws.Cells.Item(11,9).Value = 0.05996655275686102
$ws.Cells.Item(11,10).Value = 0.05996655275686102
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 1.257727
$ws.Cells.Item(11,14).Value = 3.773181
$ws.Cells.Item(11,15).Value = 0.3633293041311343
$ws.Cells.Item(11,16).Value = 0.3633293041311343
$ws.Cells.Item(11,17).Value = 31.384310860946
$ws.Cells.Item(11,18).Value = 282.458797748514
$ws.Cells.Item(11,19).Value = 0.02178760588429327
$ws.Cells.Item(11,20).Value = 0.02178760588429327

# Row 12: sCs -> M2
$ws.Cells.Item(12,1).Value = "sCs"
$ws.Cells.Item(12,2).Value = "App"
$ws.Cells.Item(12,3).Value = "Fpr2"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 24.953198
$ws.Cells.Item(12,8).Value = 74.859594
$ws.Cells.Item(12,9).Value = 0.05996655275686102
$ws.Cells.Item(12,10).Value = 0.05996655275686102
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 2.112352666666667
$ws.Cells.Item(12,14).Value = 6.337058
$ws.Cells.Item(12,15).Value = 0.6102116154455982
$ws.Cells.Item(12,16).Value = 0.6102116154455982
$ws.Cells.Item(12,17).Value = 52.70995433716134
$ws.Cells.Item(12,18).Value = 474.389589034452
$ws.Cells.Item(12,19).Value = 0.03659228703046786
$ws.Cells.Item(12,20).Value = 0.03659228703046785

# Row 13: sCs -> sCs
$ws.Cells.Item(13,1).Value = "sCs"
$ws.Cells.Item(13,2).Value = "App"
$ws.Cells.Item(13,3).Value = "Fpr2"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 24.953198
$ws.Cells.Item(13,8).Value = 74.859594
$ws.Cells.Item(13,9).Value = 0.05996655275686102
$ws.Cells.Item(13,10).Value = 0.05996655275686102
$ws.Cells.Item(13,11).Value = 1
$ws.Cells.Item(13,12).Value = 0.3333333333333333
$ws.Cells.Item(13,13).Value = 0.09159266666666667
$ws.Cells.Item(13,14).Value = 0.274778
$ws.Cells.Item(13,15).Value = 0.02645908042326748
$ws.Cells.Item(13,16).Value = 0.02645908042326749
$ws.Cells.Item(13,17).Value = 2.285529946681333
$ws.Cells.Item(13,18).Value = 20.569769520132
$ws.Cells.Item(13,19).Value = 0.001586659842099898
$ws.Cells.Item(13,20).Value = 0.001586659842099898
